$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.90"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.10"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.412"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05922"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.455"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.573"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8132"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9085"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1403"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07419"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03299"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03055"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09343"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.868"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001578"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04668"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006093"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.004980"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009818"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.00009005"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.605"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.136"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.01119"
$ws.Range("E24").Value = "23OneONEBestin24h"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002902"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04010"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006216"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1073"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003002"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008206"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005249"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.8163"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
